$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.161.78"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "'1.902.75"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  -0.57%  "

$ws.Range("D5").Value = "'253.51"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").Value = "'0.705"
$ws.Range("E6").Value = "  +2.12%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").Value = "'41.91"
$ws.Range("E8").Value = "  +3.59%  "

$ws.Range("E9").Value = "  +1.92%  "

$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").Value = "'0.0763"
$ws.Range("E11").Value = "  +5.69%  "

$ws.Range("D12").Value = "'0.0978"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").Value = "'13.17"
$ws.Range("E13").Value = "  +4.67%  "

$ws.Range("D14").Value = "'2.179.11"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "'0.735"
$ws.Range("E15").Value = "  +3.72%  "

$ws.Range("E16").Value = "  +3.79%  "

$ws.Range("D17").Value = "'1.915.99"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "'35.149.77"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").Value = "'73.84"
$ws.Range("E19").Value = "  +1.88%  "

$ws.Range("D20").Value = "'0.0₃0845"
$ws.Range("E20").Value = "  +2.91%  "

$ws.Range("D21").Value = "'243.58"
$ws.Range("E21").Value = "  +1.18%  "

$ws.Range("D22").Value = "'13.08"
$ws.Range("E22").Value = "  +2.77%  "

$ws.Range("D23").Value = "'5.05"
$ws.Range("E23").Value = "  +4.67%  "

$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").Value = "'2.43"
$ws.Range("E25").Value = "  +4.56%  "

$ws.Range("E26").Value = "  -0.76%  "

$ws.Range("D27").Value = "'168.98"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "'18.53"
$ws.Range("E29").Value = "  -2.29%  "

$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("D32").Value = "'2.11"
$ws.Range("E32").Value = "  +12.93%  "

$ws.Range("E33").Value = "  +4.26%  "

$ws.Range("D34").Value = "'0.0598"
$ws.Range("E34").Value = "  +5.18%  "

$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +6.98%  "

$ws.Range("D36").Value = "'4.26"
$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("D38").Value = "'0.846"
$ws.Range("E38").Value = "  -7.44%  "

$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'99.02"
$ws.Range("E40").Value = "  +5.06%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.32"
$ws.Range("E41").Value = "  +5.80%  "

$ws.Range("E42").Value = "  +4.09%  "

$ws.Range("D43").Value = "'0.0672"
$ws.Range("E43").Value = "  +2.99%  "

$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").Value = "'1.306.77"
$ws.Range("E46").Value = "  -3.54%  "

$ws.Range("D47").Value = "'2.41"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").Value = "'2.74"
$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("D49").Value = "'6.60"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").Value = "'12.14"
$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("D51").Value = "'0.0759"
$ws.Range("E51").Value = "  +7.95%  "
